# LOM3223.xlsx restructuring:
# The "Programa" / "Avaliação" block (old rows 13-24) is reshuffled: a new
# column-A label is introduced on rows that previously had none, some
# B/C value pairs are dropped, others are added, row heights change, and
# the two trailing rows (25-26) are removed entirely (their content is
# folded into what is now row 23-24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Drop the two now-unused trailing rows first (this also re-derives
#    the sheet dimension to A1:C24 automatically).
# ---------------------------------------------------------------------
$ws.Range("A26:C26").EntireRow.Delete() | Out-Null
$ws.Range("A25:C25").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2) Row 13: gains a label in column A ("Programa resumido:"); B/C keep
#    their existing formatting but now show the activation date.
# ---------------------------------------------------------------------
$ws.Range("A13").Style = $ws.Range("A12").Style
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2012"
$ws.Range("C13").Value = "01/01/2012"

# ---------------------------------------------------------------------
# 3) Row 14: gains a label in column A ("Short syllabus:"); old B/C pair
#    is removed entirely.
# ---------------------------------------------------------------------
$ws.Range("A14").Style = $ws.Range("A12").Style
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear() | Out-Null
$ws.Range("C14").Clear() | Out-Null

# ---------------------------------------------------------------------
# 4) Row 15: label stays ("Programa:"), B/C now show the responsible
#    professor instead of the long syllabus text.
# ---------------------------------------------------------------------
$ws.Range("B15").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C15").Value = "5840730 - Antonio Jefferson da Silva Machado"

# ---------------------------------------------------------------------
# 5) Row 16: label text changes to "Syllabus:".
# ---------------------------------------------------------------------
$ws.Range("A16").Value = "Syllabus:"

# ---------------------------------------------------------------------
# 6) Row 17: label text changes to "Avaliação:"; B/C pair removed.
# ---------------------------------------------------------------------
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear() | Out-Null
$ws.Range("C17").Clear() | Out-Null

# ---------------------------------------------------------------------
# 7) Row 18: label text changes to "Método:"; gains a new B/C pair with
#    the second professor's name.
# ---------------------------------------------------------------------
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Style = $ws.Range("B20").Style
$ws.Range("C18").Style = $ws.Range("C20").Style
$ws.Range("B18").Value = "5840726 - Cristina Bormio Nunes"
$ws.Range("C18").Value = "5840726 - Cristina Bormio Nunes"

# ---------------------------------------------------------------------
# 8) Row 19: label text changes to "Critério:"; gains a new B/C pair
#    with the teaching-method description.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Style = $ws.Range("B20").Style
$ws.Range("C19").Style = $ws.Range("C20").Style
$ws.Range("B19").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C19").Value = "Aulas expositivas, seminários e exercícios comentados."

# Rows 20 and 21 keep their text content (only row 21's height changes,
# handled below).

# ---------------------------------------------------------------------
# 9) Row 22: B/C pair removed, label text unchanged ("Requisitos:" stays
#    on what is now row 22... actually label unchanged = "Norma de
#    recuperação:" remains as-is here).
# ---------------------------------------------------------------------
$ws.Range("B22").Clear() | Out-Null
$ws.Range("C22").Clear() | Out-Null

# ---------------------------------------------------------------------
# 10) Row 23: column-A label removed; B/C now hold the first requirement
#     line instead of the bibliography text.
# ---------------------------------------------------------------------
$ws.Range("A23").Clear() | Out-Null
$ws.Range("B23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Range("C23").Value = "LOM3206 -  Eletrônica  (Requisito)`n"

# ---------------------------------------------------------------------
# 11) Row 24: column-A label removed; gains a new B/C pair with the
#     second requirement line.
# ---------------------------------------------------------------------
$ws.Range("A24").Clear() | Out-Null
$ws.Range("B24").Style = $ws.Range("B23").Style
$ws.Range("C24").Style = $ws.Range("C23").Style
$ws.Range("B24").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$ws.Range("C24").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"

# ---------------------------------------------------------------------
# 12) Row heights per the new layout.
# ---------------------------------------------------------------------
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(17).AutoFit() | Out-Null
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(22).AutoFit() | Out-Null
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30
